$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1861
$ws.Range("J3").Value = 8078
$ws.Range("K3").Value = 1774
$ws.Range("J4").Value = 1801
$ws.Range("K4").Value = 386
$ws.Range("K5").Value = 118
$ws.Range("J6").Value = 11059
$ws.Range("K6").Value = 2288
$ws.Range("J7").Value = 29270
$ws.Range("K7").Value = 6427

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 178
$ws.Range("K8").Value = 420
$ws.Range("K13").Value = 10
$ws.Range("K14").Value = 36
$ws.Range("K17").Value = 10
$ws.Range("K18").Value = 47
$ws.Range("K19").Value = 177
$ws.Range("J20").Value = 636
$ws.Range("K20").Value = 142
$ws.Range("K22").Value = 17
$ws.Range("K23").Value = 58
$ws.Range("K29").Value = 310
$ws.Range("K31").Value = 73
$ws.Range("K33").Value = 260
$ws.Range("K35").Value = 10
$ws.Range("K37").Value = 212
$ws.Range("K42").Value = 221
$ws.Range("K43").Value = 62
$ws.Range("K44").Value = 62
$ws.Range("K46").Value = 13
$ws.Range("J48").Value = 322
$ws.Range("K48").Value = 77
$ws.Range("K49").Value = 38
$ws.Range("K50").Value = 35
$ws.Range("K51").Value = 75
$ws.Range("K53").Value = 101
$ws.Range("J54").Value = 574
$ws.Range("K54").Value = 107
$ws.Range("K55").Value = 68
$ws.Range("K56").Value = 10
$ws.Range("K57").Value = 14
$ws.Range("K60").Value = 43
$ws.Range("K64").Value = 44
$ws.Range("K65").Value = 156
$ws.Range("K67").Value = 248
$ws.Range("K72").Value = 26
$ws.Range("K75").Value = 23
$ws.Range("K76").Value = 92
$ws.Range("K79").Value = 172
$ws.Range("K83").Value = 135
$ws.Range("K85").Value = 324
$ws.Range("K86").Value = 44
$ws.Range("K91").Value = 57
$ws.Range("K94").Value = 76
$ws.Range("K96").Value = 92
$ws.Range("K97").Value = 58
$ws.Range("J101").Value = 29270
$ws.Range("K101").Value = 6427

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 64
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 116
$ws.Range("K3").Value = 107
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 21
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 122
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 420

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 98
$ws.Range("K7").Value = 260

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 37
$ws.Range("J6").Value = 263
$ws.Range("J7").Value = 574
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 83
$ws.Range("K3").Value = 106
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 310

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 12
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 14
$ws.Range("K6").Value = 35
$ws.Range("J7").Value = 322
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 221

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 10

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 60
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 10
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 42
$ws.Range("J3").Value = 209
$ws.Range("K6").Value = 52
$ws.Range("J7").Value = 636
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 10
